$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 557.06665
$ws.Range("I33").Value = 129.66667
$ws.Range("J33").Value = 2266.6667
$ws.Range("K33").Value = 129.66667
$ws.Range("L33").Value = 2266.6667
$ws.Range("M33").Value = 99.33332999999999
$ws.Range("N33").Value = -2724.6667
$ws.Range("H40").Value = 2750
$ws.Range("J40").Value = 3000
$ws.Range("L40").Value = 3000
$ws.Range("N40").Value = -3350
$ws.Range("H61").Value = 1698.75
$ws.Range("I61").Value = 1698.75
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 5096.25
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -4924.25
$ws.Range("N61").ClearContents()
$ws.Range("H86").Value = 4248.25
$ws.Range("I86").Value = 4759.4
$ws.Range("J86").Value = 3883.1428
$ws.Range("K86").Value = 4759.4
$ws.Range("L86").Value = 3883.1428
$ws.Range("M86").Value = -3636.4
$ws.Range("N86").Value = -6129.1428
$ws.Range("H89").Value = 4248.25
$ws.Range("I89").Value = 4759.4
$ws.Range("J89").Value = 3883.1428
$ws.Range("K89").Value = 23797
$ws.Range("L89").Value = 19415.714
$ws.Range("M89").Value = -18181
$ws.Range("N89").Value = -30647.714
$ws.Range("H100").Value = 1992.5
$ws.Range("I100").Value = 1992.5
$ws.Range("K100").Value = 1992.5
$ws.Range("M100").Value = -1451.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2795.3225
$ws.Range("I32").Value = 2348.3928
$ws.Range("K32").Value = 2348.3928
$ws.Range("M32").Value = -2061.3928
$ws.Range("H46").Value = 28381.8
$ws.Range("J46").Value = 32636.334
$ws.Range("L46").Value = 32636.334
$ws.Range("N46").Value = -33274.334
$ws.Range("H61").Value = 5000
$ws.Range("J61").Value = 5000
$ws.Range("L61").Value = 5000
$ws.Range("N61").Value = -5424
$ws.Range("H92").Value = 67997.25
$ws.Range("J92").Value = 67997.25
$ws.Range("L92").Value = 67997.25
$ws.Range("N92").Value = -72989.25
$ws.Range("H132").Value = 3900
$ws.Range("I132").Value = 3900
$ws.Range("K132").Value = 11700
$ws.Range("M132").Value = -9170
$ws.Range("H136").Value = 5000
$ws.Range("J136").Value = 5000
$ws.Range("L136").Value = 15000
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 876.2222
$ws.Range("I20").Value = 798.75
$ws.Range("K20").Value = 798.75
$ws.Range("M20").Value = -551.75
$ws.Range("H26").Value = 7900
$ws.Range("I26").Value = 7900
$ws.Range("K26").Value = 7900
$ws.Range("M26").Value = -7608
$ws.Range("H62").Value = 45181
$ws.Range("J62").Value = 45181
$ws.Range("L62").Value = 45181
$ws.Range("N62").Value = -46553
$ws.Range("H65").Value = 45181
$ws.Range("J65").Value = 45181
$ws.Range("L65").Value = 135543
$ws.Range("N65").Value = -142407
$ws.Range("H94").Value = 4169.3335
$ws.Range("I94").Value = 3754.5
$ws.Range("K94").Value = 3754.5
$ws.Range("M94").Value = -3303.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1811.037
$ws.Range("I31").Value = 1616.2632
$ws.Range("J31").Value = 2273.625
$ws.Range("K31").Value = 1616.2632
$ws.Range("L31").Value = 2273.625
$ws.Range("M31").Value = -1321.2632
$ws.Range("N31").Value = -2863.625
$ws.Range("H34").Value = 1811.037
$ws.Range("I34").Value = 1616.2632
$ws.Range("J34").Value = 2273.625
$ws.Range("K34").Value = 1616.2632
$ws.Range("L34").Value = 2273.625
$ws.Range("M34").Value = -1414.2632
$ws.Range("N34").Value = -2677.625
$ws.Range("H36").Value = 3950
$ws.Range("J36").Value = 3950
$ws.Range("L36").Value = 3950
$ws.Range("N36").Value = -4726
$ws.Range("H40").Value = 3950
$ws.Range("J40").Value = 3950
$ws.Range("L40").Value = 3950
$ws.Range("N40").Value = -4270
$ws.Range("H51").Value = 16333.333
$ws.Range("J51").Value = 20000
$ws.Range("L51").Value = 20000
$ws.Range("N51").Value = -21472
$ws.Range("H54").Value = 9999
$ws.Range("J54").Value = 9999
$ws.Range("L54").Value = 9999
$ws.Range("N54").Value = -11315
$ws.Range("H58").Value = 3295.2856
$ws.Range("I58").Value = 3511.3333
$ws.Range("K58").Value = 3511.3333
$ws.Range("M58").Value = -3308.3333
$ws.Range("H61").Value = 16333.333
$ws.Range("J61").Value = 20000
$ws.Range("L61").Value = 20000
$ws.Range("N61").Value = -20696
$ws.Range("H136").Value = 3295.2856
$ws.Range("I136").Value = 3511.3333
$ws.Range("K136").Value = 10533.9999
$ws.Range("M136").Value = -7983.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1500
$ws.Range("J86").Value = 1500
$ws.Range("L86").Value = 4500
$ws.Range("N86").Value = -6872
$ws.Range("H89").Value = 1500
$ws.Range("J89").Value = 1500
$ws.Range("L89").Value = 13500
$ws.Range("N89").Value = -25356
$ws.Range("H109").Value = 3485.7144
$ws.Range("J109").Value = 3485.7144
$ws.Range("L109").Value = 10457.1432
$ws.Range("N109").Value = -12537.1432
$ws.Range("H129").Value = 1438.4
$ws.Range("I129").Value = 772.3333
$ws.Range("J129").Value = 2437.5
$ws.Range("K129").Value = 2316.9999
$ws.Range("L129").Value = 7312.5
$ws.Range("M129").Value = 2683.0001
$ws.Range("N129").Value = -17312.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 100.47059
$ws.Range("J2").Value = 145.22223
$ws.Range("L2").Value = 145.22223
$ws.Range("N2").Value = -371.22223
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("N24").ClearContents()
$ws.Range("H92").Value = 13500
$ws.Range("J92").Value = 13500
$ws.Range("L92").Value = 13500
$ws.Range("N92").Value = -17244
$ws.Range("H97").Value = 3271
$ws.Range("I97").Value = 3210.5
$ws.Range("J97").Value = 3432.3333
$ws.Range("K97").Value = 3210.5
$ws.Range("L97").Value = 3432.3333
$ws.Range("M97").Value = -2714.5
$ws.Range("N97").Value = -4424.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 20000
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H5").Value = 18500
$ws.Range("I5").Value = 35000
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 35000
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = -34887
$ws.Range("N5").Value = -2226
$ws.Range("H15").Value = 20000
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H93").Value = 2967.6667
$ws.Range("I93").Value = 2967.6667
$ws.Range("K93").Value = 2967.6667
$ws.Range("M93").Value = -1719.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H30").Value = 5000
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 5000
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 5000
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -5214
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H92").Value = 29800
$ws.Range("J92").Value = 29800
$ws.Range("L92").Value = 29800
$ws.Range("N92").Value = -34792
